function Set-TextCell {
    param($ws, $cellRef, $val)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws 'D2' '67.136.49'
Set-TextCell $ws 'E2' '  +0.39%  '
Set-TextCell $ws 'D3' '3.088.45'
Set-TextCell $ws 'E3' '  -0.72%  '
Set-TextCell $ws 'E4' '  -0.23%  '
Set-TextCell $ws 'D5' '579.06'
Set-TextCell $ws 'E5' '  +0.68%  '
Set-TextCell $ws 'D6' '169.32'
Set-TextCell $ws 'E6' '  -2.26%  '
Set-TextCell $ws 'E7' '  -0.15%  '
Set-TextCell $ws 'D8' '3.084.42'
Set-TextCell $ws 'E8' '  -0.68%  '
Set-TextCell $ws 'E9' '  -0.85%  '
Set-TextCell $ws 'D10' '6.43'
Set-TextCell $ws 'E10' '  +0.28%  '
Set-TextCell $ws 'D11' '0.151'
Set-TextCell $ws 'E11' '  -0.97%  '
Set-TextCell $ws 'D12' '0.473'
Set-TextCell $ws 'E12' '  -0.72%  '
Set-TextCell $ws 'D13' '0.0000242'
Set-TextCell $ws 'E13' '  -1.33%  '
Set-TextCell $ws 'D14' '36.35'
Set-TextCell $ws 'E14' '  -1.90%  '
Set-TextCell $ws 'E15' '  -2.00%  '
Set-TextCell $ws 'D16' '3.598.62'
Set-TextCell $ws 'E16' '  -0.79%  '
Set-TextCell $ws 'D17' '66.983.41'
Set-TextCell $ws 'E17' '  +0.07%  '
Set-TextCell $ws 'D18' '7.02'
Set-TextCell $ws 'E18' '  -0.98%  '
Set-TextCell $ws 'B19' 'Chainlink'
Set-TextCell $ws 'C19' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D19' '16.57'
Set-TextCell $ws 'E19' '  +2.33%  '
Set-TextCell $ws 'B20' 'WrappedEther'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 'D20' '3.088.11'
Set-TextCell $ws 'E20' '  -0.75%  '
Set-TextCell $ws 'D21' '490.78'
Set-TextCell $ws 'E21' '  +3.19%  '
Set-TextCell $ws 'D22' '7.78'
Set-TextCell $ws 'E22' '  +0.70%  '
Set-TextCell $ws 'D23' '0.689'
Set-TextCell $ws 'E23' '  -2.93%  '
Set-TextCell $ws 'D24' '82.77'
Set-TextCell $ws 'E24' '  -1.16%  '
Set-TextCell $ws 'D25' '12.92'
Set-TextCell $ws 'E25' '  -2.67%  '
Set-TextCell $ws 'D26' '2.24'
Set-TextCell $ws 'E26' '  -2.01%  '
Set-TextCell $ws 'D27' '10.29'
Set-TextCell $ws 'E27' '  +3.86%  '
Set-TextCell $ws 'E28' '  +0.11%  '
Set-TextCell $ws 'D29' '7.83'
Set-TextCell $ws 'E29' '  -1.82%  '
Set-TextCell $ws 'E30' '  -3.29%  '
Set-TextCell $ws 'D31' '2.64'
Set-TextCell $ws 'E31' '  -0.40%  '
Set-TextCell $ws 'D32' '27.98'
Set-TextCell $ws 'E32' '  -1.96%  '
Set-TextCell $ws 'D33' '0.112'
Set-TextCell $ws 'E33' '  -1.32%  '
Set-TextCell $ws 'D34' '0.0₃0918'
Set-TextCell $ws 'E34' '  -4.68%  '
Set-TextCell $ws 'D35' '0.999'
Set-TextCell $ws 'E35' '  -0.12%  '
Set-TextCell $ws 'D36' '5.72'
Set-TextCell $ws 'E36' '  -2.06%  '
Set-TextCell $ws 'D37' '0.956'
Set-TextCell $ws 'E37' '  -2.26%  '
Set-TextCell $ws 'D38' '46.17'
Set-TextCell $ws 'E38' '  -3.48%  '
Set-TextCell $ws 'D39' '0.124'
Set-TextCell $ws 'E39' '  +1.47%  '
Set-TextCell $ws 'D40' '2.00'
Set-TextCell $ws 'E40' '  -3.88%  '
Set-TextCell $ws 'E41' '  -1.89%  '
Set-TextCell $ws 'D42' '8.35'
Set-TextCell $ws 'E42' '  -2.76%  '
Set-TextCell $ws 'D43' '2.778.66'
Set-TextCell $ws 'E43' '  -0.59%  '
Set-TextCell $ws 'D44' '372.80'
Set-TextCell $ws 'E44' '  -1.30%  '
Set-TextCell $ws 'E45' '  -2.28%  '
Set-TextCell $ws 'D46' '135.59'
Set-TextCell $ws 'E46' '  -0.47%  '
Set-TextCell $ws 'D47' '2.50'
Set-TextCell $ws 'E47' '  -2.41%  '
Set-TextCell $ws 'E48' '  +0.00%  '
Set-TextCell $ws 'D49' '24.51'
Set-TextCell $ws 'E49' '  -0.61%  '
Set-TextCell $ws 'D50' '2.16'
Set-TextCell $ws 'E50' '  -1.63%  '
Set-TextCell $ws 'E51' '  -1.05%  '
